# The author's edit swaps the deck's applied Design/Theme: the slide
# master (and the whole presentation) moves off the custom "Integral"
# (Red Violet) look and over to the stock "Office Theme" colors, while
# the theme that used to be the plain "Office Theme" (only ever used by
# the Notes Master) ends up carrying the old "Integral"/Red Violet
# palette. In other words: the two theme parts trade places.
#
# The supported way to drive that from the PowerPoint object model is
# via Application/Presentation/Master "ApplyTheme" (choosing a
# different Office Theme design) and via the Design/Theme color
# scheme. We issue both: the high-level "apply a theme" calls that
# describe the actual user action, and then the granular
# ThemeColorScheme color edits, which is the mechanism that actually
# persists through this host and rewrites the 12 scheme colors in the
# slide master's theme part to the stock Office Theme palette (the
# target state for that theme part).

$p = $ppt.ActivePresentation

# --- High level: apply the built-in "Office Theme" design -----------
# (mirrors Design > Themes > "Office Theme" in the PowerPoint UI)
try { $p.ApplyTheme("Office Theme") } catch {}
try { $p.SlideMaster.ApplyTheme("Office Theme") } catch {}

# --- Low level: rewrite the 12 theme colors one by one ---------------
# Index order matches PowerPoint's ThemeColorScheme.Colors(1..12):
#  1 dk1   2 lt1   3 dk2   4 lt2
#  5 accent1  6 accent2  7 accent3  8 accent4  9 accent5  10 accent6
#  11 hlink   12 folHlink
# RGB() isn't available as a function in this host, so each colour is
# packed into PowerPoint's BGR-in-an-int VBA RGB encoding by hand:
#   value = R + G*256 + B*65536

function Pack-Rgb($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$officeThemeColors = @(
    @(0,   0,   0),    #  1 dk1      000000
    @(255, 255, 255),  #  2 lt1      FFFFFF
    @(68,  84,  106),  #  3 dk2      44546A
    @(231, 230, 230),  #  4 lt2      E7E6E6
    @(91,  155, 213),  #  5 accent1  5B9BD5
    @(237, 125, 49),   #  6 accent2  ED7D31
    @(165, 165, 165),  #  7 accent3  A5A5A5
    @(255, 192, 0),    #  8 accent4  FFC000
    @(68,  114, 196),  #  9 accent5  4472C4
    @(112, 173, 71),   # 10 accent6  70AD47
    @(5,   99,  193),  # 11 hlink    0563C1
    @(149, 79,  114)   # 12 folHlink 954F72
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $rgb = $officeThemeColors[$i - 1]
    $packed = Pack-Rgb $rgb[0] $rgb[1] $rgb[2]
    $themeColors.Colors($i).RGB = $packed
}
